$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 2 (old row2->4, old row3->5),
# pushing the existing "Code Snippet ID"/"Link" header row and the blank
# data row down to make room for two new label rows.
$ws.Rows("2:3").Insert()

# The insert copies formatting from the row above into the new rows'
# interior columns (C2:C3) - clear that stray formatting before writing
# the real content.
$ws.Range("C2:C3").Clear()

# New row 2 / row 3 content: "Suspicious" / "Critical" labels in column B,
# styled like the existing "Metrics" header (centered, same row height).
$ws.Range("B2").Value = "Suspicious"
$ws.Range("B3").Value = "Critical"
$ws.Range("B2:B3").HorizontalAlignment = -4108   # xlCenter
$ws.Rows("2:3").RowHeight = 15.75

# Re-point the existing conditional formats at their shifted ranges
# (A3:D3 E2:F2 -> A5:D5 E4:F4, H2 -> H4). The COM bridge's
# ModifyAppliesToRange only keeps the first contiguous area of a
# multi-area range, so re-target the first area on the original rule and
# add a twin rule (same test/fill) for the second area.
$fc1 = $ws.Cells.FormatConditions.Item(1)
$fc1.ModifyAppliesToRange($ws.Range("A5:D5"))

$fc2 = $ws.Cells.FormatConditions.Item(2)
$fc2.ModifyAppliesToRange($ws.Range("H4"))

$fcExtra = $ws.Range("E4:F4").FormatConditions.Add(1, 3, '"Yes"')
$fcExtra.Interior.Color = 6711008

# Refresh the frozen-pane split and selection to match the new layout
# (freeze through row 4, active cell at D12).
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D12").Select()
